$p = $ppt.ActivePresentation

# ppPlaceholderDate = 16
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes, [string]$newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch {
        }
        if ($isDate -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# Update the date placeholder text on the slide master ...
Update-DatePlaceholder $p.SlideMaster.Shapes "12/31/09"

# ... and on every slide layout (custom layout) tied to the master.
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes "12/31/09"
}

# Update the "ACE + TAO" cube label to "DOC" on the first slide.
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $sh = $slide1.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "ACE + TAO") {
            $sh.TextFrame.TextRange.Text = "DOC"
        }
    }
}
